# Auto-generated: apply scheduled market-price refresh to Leve profit tables
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 997.5
$ws.Range("I19").Value = 995
$ws.Range("K19").Value = 995
$ws.Range("M19").Value = -820
$ws.Range("H28").Value = 517
$ws.Range("I28").Value = 526.6667
$ws.Range("K28").Value = 526.6667
$ws.Range("M28").Value = -41.66669999999999
$ws.Range("H53").Value = 305.4
$ws.Range("I53").Value = 221.33333
$ws.Range("K53").Value = 221.33333
$ws.Range("M53").Value = 415.66667
$ws.Range("H98").Value = 3335.9092
$ws.Range("I98").Value = 1711.8889
$ws.Range("J98").Value = 10644
$ws.Range("K98").Value = 1711.8889
$ws.Range("L98").Value = 10644
$ws.Range("M98").Value = -213.8888999999999
$ws.Range("N98").Value = -13640
$ws.Range("H107").Value = 2026.6154
$ws.Range("I107").Value = 2026.6154
$ws.Range("K107").Value = 2026.6154
$ws.Range("M107").Value = -106.6153999999999
$ws.Range("H112").Value = 3061.111
$ws.Range("J112").Value = 3318.75
$ws.Range("L112").Value = 9956.25
$ws.Range("N112").Value = -12172.25
$ws.Range("H113").Value = 3606.625
$ws.Range("I113").Value = 3100.75
$ws.Range("J113").Value = 4112.5
$ws.Range("K113").Value = 3100.75
$ws.Range("L113").Value = 4112.5
$ws.Range("M113").Value = 153.25
$ws.Range("N113").Value = -10620.5
$ws.Range("H122").Value = 3335.9092
$ws.Range("I122").Value = 1711.8889
$ws.Range("J122").Value = 10644
$ws.Range("K122").Value = 5135.6667
$ws.Range("L122").Value = 31932
$ws.Range("M122").Value = -2685.6667
$ws.Range("N122").Value = -36832
$ws.Range("H132").Value = 1598.069
$ws.Range("I132").Value = 1598.069
$ws.Range("K132").Value = 4794.207
$ws.Range("M132").Value = -2264.207
$ws.Range("H137").Value = 2455.3
$ws.Range("I137").Value = 2375.6667
$ws.Range("J137").Value = 2574.75
$ws.Range("K137").Value = 7127.000100000001
$ws.Range("L137").Value = 7724.25
$ws.Range("M137").Value = -4577.000100000001
$ws.Range("N137").Value = -12824.25
$ws.Range("H138").Value = 7842
$ws.Range("J138").Value = 8009.1763
$ws.Range("L138").Value = 24027.5289
$ws.Range("N138").Value = -34307.5289

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("H32").Value = 11077.098
$ws.Range("I32").Value = 8833.172
$ws.Range("J32").Value = 24166.666
$ws.Range("K32").Value = 8833.172
$ws.Range("L32").Value = 24166.666
$ws.Range("M32").Value = -8546.172
$ws.Range("N32").Value = -24740.666
$ws.Range("H45").Value = 2960.5
$ws.Range("I45").Value = 2960.5
$ws.Range("K45").Value = 2960.5
$ws.Range("M45").Value = -2583.5
$ws.Range("H46").Value = 14822.25
$ws.Range("J46").Value = 19575.5
$ws.Range("L46").Value = 19575.5
$ws.Range("N46").Value = -20213.5
$ws.Range("H61").Value = 3832.3333
$ws.Range("I61").Value = 3832.3333
$ws.Range("K61").Value = 3832.3333
$ws.Range("M61").Value = -3620.3333
$ws.Range("H74").Value = 6301.8335
$ws.Range("I74").Value = 7162.2
$ws.Range("K74").Value = 7162.2
$ws.Range("M74").Value = -6288.2
$ws.Range("H77").Value = 6301.8335
$ws.Range("I77").Value = 7162.2
$ws.Range("K77").Value = 35811
$ws.Range("M77").Value = -31443
$ws.Range("H132").Value = 1598.0667
$ws.Range("I132").Value = 747.9167
$ws.Range("K132").Value = 2243.7501
$ws.Range("M132").Value = 286.2498999999998
$ws.Range("H136").Value = 3832.3333
$ws.Range("I136").Value = 3832.3333
$ws.Range("K136").Value = 11496.9999
$ws.Range("M136").Value = -8946.999899999999

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4903.5557
$ws.Range("I94").Value = 5890.3335
$ws.Range("J94").Value = 2930
$ws.Range("K94").Value = 5890.3335
$ws.Range("L94").Value = 2930
$ws.Range("M94").Value = -5439.3335
$ws.Range("N94").Value = -3832
$ws.Range("H99").Value = 35336.668
$ws.Range("I99").Value = 35336.668
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 35336.668
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -33838.668
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 2175.8333
$ws.Range("I105").Value = 2175.8333
$ws.Range("K105").Value = 2175.8333
$ws.Range("M105").Value = -428.8332999999998
$ws.Range("H107").Value = 2236.7
$ws.Range("I107").Value = 2150.2856
$ws.Range("J107").Value = 2438.3333
$ws.Range("K107").Value = 2150.2856
$ws.Range("L107").Value = 2438.3333
$ws.Range("M107").Value = -230.2856000000002
$ws.Range("N107").Value = -6278.3333
$ws.Range("H134").Value = 1472.625
$ws.Range("I134").Value = 1289.3846
$ws.Range("K134").Value = 3868.1538
$ws.Range("M134").Value = -1333.1538

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2354.5
$ws.Range("I31").Value = 2081.8333
$ws.Range("J31").Value = 3172.5
$ws.Range("K31").Value = 2081.8333
$ws.Range("L31").Value = 3172.5
$ws.Range("M31").Value = -1786.8333
$ws.Range("N31").Value = -3762.5
$ws.Range("H34").Value = 2354.5
$ws.Range("I34").Value = 2081.8333
$ws.Range("J34").Value = 3172.5
$ws.Range("K34").Value = 2081.8333
$ws.Range("L34").Value = 3172.5
$ws.Range("M34").Value = -1879.8333
$ws.Range("N34").Value = -3576.5
$ws.Range("H58").Value = 3444.625
$ws.Range("I58").Value = 3411.4
$ws.Range("K58").Value = 3411.4
$ws.Range("M58").Value = -3208.4
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H74").Value = 71314
$ws.Range("J74").Value = 71314
$ws.Range("L74").Value = 71314
$ws.Range("N74").Value = -73062
$ws.Range("H77").Value = 71314
$ws.Range("J77").Value = 71314
$ws.Range("L77").Value = 213942
$ws.Range("N77").Value = -222678
$ws.Range("H122").Value = 2403.7273
$ws.Range("I122").Value = 2348.8235
$ws.Range("K122").Value = 7046.470499999999
$ws.Range("M122").Value = -4596.470499999999
$ws.Range("H134").Value = 2168.4348
$ws.Range("I134").Value = 1666.6111
$ws.Range("K134").Value = 4999.8333
$ws.Range("M134").Value = -2464.8333
$ws.Range("H136").Value = 3444.625
$ws.Range("I136").Value = 3411.4
$ws.Range("K136").Value = 10234.2
$ws.Range("M136").Value = -7684.200000000001

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 480.16666
$ws.Range("I38").Value = 465.33334
$ws.Range("J38").Value = 495
$ws.Range("K38").Value = 1396.00002
$ws.Range("L38").Value = 1485
$ws.Range("M38").Value = -1049.00002
$ws.Range("N38").Value = -2179
$ws.Range("H98").Value = 2646.4285
$ws.Range("I98").Value = 2859.1667
$ws.Range("J98").Value = 1370
$ws.Range("K98").Value = 8577.500100000001
$ws.Range("L98").Value = 4110
$ws.Range("M98").Value = -7079.500100000001
$ws.Range("N98").Value = -7106

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3996.25
$ws.Range("I80").Value = 3997
$ws.Range("K80").Value = 3997
$ws.Range("M80").Value = -2999
$ws.Range("H83").Value = 3996.25
$ws.Range("I83").Value = 3997
$ws.Range("K83").Value = 19985
$ws.Range("M83").Value = -14993
$ws.Range("H126").Value = 2723.7144
$ws.Range("I126").Value = 2943.7
$ws.Range("J126").Value = 2173.75
$ws.Range("K126").Value = 8831.099999999999
$ws.Range("L126").Value = 6521.25
$ws.Range("M126").Value = -6361.099999999999
$ws.Range("N126").Value = -11461.25
$ws.Range("H132").Value = 2711.5454
$ws.Range("I132").Value = 1876.1428
$ws.Range("K132").Value = 5628.428400000001
$ws.Range("M132").Value = -3098.428400000001
$ws.Range("H141").Value = 98998.336
$ws.Range("J141").Value = 98998.336
$ws.Range("L141").Value = 98998.336
$ws.Range("N141").Value = -109358.336

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8039.5713
$ws.Range("I7").Value = 4379.6665
$ws.Range("K7").Value = 4379.6665
$ws.Range("M7").Value = -4267.6665
$ws.Range("H46").Value = 200
$ws.Range("I46").Value = 200
$ws.Range("K46").Value = 200
$ws.Range("M46").Value = -12
$ws.Range("H126").Value = 8039.5713
$ws.Range("I126").Value = 4379.6665
$ws.Range("K126").Value = 13138.9995
$ws.Range("M126").Value = -10668.9995

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1278.5714
$ws.Range("I100").Value = 1369
$ws.Range("J100").Value = 1158
$ws.Range("K100").Value = 2738
$ws.Range("L100").Value = 2316
$ws.Range("M100").Value = -2197
$ws.Range("N100").Value = -3398
$ws.Range("H132").Value = 3592
$ws.Range("I132").Value = 1889.7693
$ws.Range("K132").Value = 5669.3079
$ws.Range("M132").Value = -3139.3079
$ws.Range("H136").Value = 983.82355
$ws.Range("I136").Value = 901.6667
$ws.Range("K136").Value = 2705.0001
$ws.Range("M136").Value = -155.0001000000002

